$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59-64 down to 60-65
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record's data
$ws.Range("A59").Value = 2
$ws.Range("B59").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44769
$ws.Range("D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E59").Value = 4
$ws.Range("F59").Value = 100112022
$ws.Range("G59").Value = "Arveja Verde"
$ws.Range("H59").Value = "Perfection"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 500
$ws.Range("K59").Value = 30000
$ws.Range("L59").Value = 32000
$ws.Range("M59").Value = 31000
$ws.Range("N59").Value = "`$/malla 25 kilos"
$ws.Range("O59").Value = "Provincia de Limarí"
$ws.Range("P59").Value = 1240
$ws.Range("Q59").Value = 25
$ws.Range("R59").Value = "Hortaliza"
